# Generate Report for Handback
# Updates the localization-status report to reflect that the de-de / zh-cn
# handback packages have now been generated and are in sync with en-US.

$wb = $excel.ActiveWorkbook

$hyperlinkTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0bf409f3215d9b69f343e04893ccb9b80ccadf42/e2e/2884f280-caf2-4e65-8a5f-da671b7c46a8.md"
$targetFileName  = "2884f280-caf2-4e65-8a5f-da671b7c46a8.md"
$statusText      = "Handed back: in sync with en-US"

# --- zh-cn sheet -----------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Status -> Handed back
$wsZh.Range("C2").Value2 = $statusText

# Latest Target File (+hyperlink to the source markdown file)
$wsZh.Range("I2").Value2 = $targetFileName
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $hyperlinkTarget, [Type]::Missing, [Type]::Missing, $targetFileName) | Out-Null

# Latest Handback File
$wsZh.Range("J2").Value2 = "2884f280-caf2-4e65-8a5f-da671b7c46a8.a25ca0b35341d4bbd284fc6a753823d81cc73fb4.zh-cn.xlf"

# Latest Handback DateTime
$wsZh.Range("K2").Value2 = "2016-09-05 11:27:13"

# Widen columns to fit the newly-populated, longer content
$wsZh.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZh.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZh.Columns.Item(10).ColumnWidth = 39.166666666666664

# --- de-de sheet -------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Status -> Handed back
$wsDe.Range("C2").Value2 = $statusText

# Latest Target File (+hyperlink to the source markdown file)
$wsDe.Range("I2").Value2 = $targetFileName
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $hyperlinkTarget, [Type]::Missing, [Type]::Missing, $targetFileName) | Out-Null

# Latest Handback File
$wsDe.Range("J2").Value2 = "2884f280-caf2-4e65-8a5f-da671b7c46a8.a25ca0b35341d4bbd284fc6a753823d81cc73fb4.de-de.xlf"

# Latest Handback DateTime
$wsDe.Range("K2").Value2 = "2016-09-05 11:27:21"

# Widen columns to fit the newly-populated, longer content
$wsDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDe.Columns.Item(10).ColumnWidth = 39.166666666666664

# --- Overview sheet ------------------------------------------------
# The Overview sheet's zh-cn/de-de columns mirror the same "Status" shared
# string used on the per-language sheets, so they flip to the new text too.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value2 = $statusText
$wsOverview.Range("F2").Value2 = $statusText
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

Write-Host "Handback report generated."
